# UC002 test-suite workbook: swap the step/expected-result text that
# belongs to the TC2 block with the text that belongs to the TC4 block.
#
# Before:
#   TC2 -> "Beneficiário Clica em cancelar diária." / "SYSTEM Apresenta a tela de Cancelar Solicitação de Diária"
#   TC3 -> "Beneficiário Clica em detalhar diária." / "SYSTEM Apresenta a tela de Detalhar Diárias"
#   TC4 -> "Beneficiário Clica em analisar prestação de contas." / "SYSTEM Apresenta a tela de Analisar Prestação de Contas"
#
# After:
#   TC2 -> "Beneficiário Clica em analisar prestação de contas." / "SYSTEM Apresenta a tela de Analisar Prestação de Contas"
#   TC3 -> unchanged ("detalhar diária")
#   TC4 -> "Beneficiário Clica em cancelar diária." / "SYSTEM Apresenta a tela de Cancelar Solicitação de Diária"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each test case "block" starts with a header row containing
# "Test Case ID: " in column A and the case id (e.g. "TC2") in column B.
# Four rows below that header row is the data row holding the step text
# in column B ("Steps") and the expected result text in column D
# ("Expected Results").
$headerRowOffsetToData = 4

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

$tc2DataRow = $null
$tc4DataRow = $null

for ($r = 1; $r -le $rowCount; $r++) {
    $caseId = $ws.Cells.Item($r, 2).Value2
    if ($caseId -eq "TC2") {
        $tc2DataRow = $r + $headerRowOffsetToData
    } elseif ($caseId -eq "TC4") {
        $tc4DataRow = $r + $headerRowOffsetToData
    }
}

if ($tc2DataRow -ne $null -and $tc4DataRow -ne $null) {
    $tc2Step = $ws.Cells.Item($tc2DataRow, 2).Value2
    $tc2Result = $ws.Cells.Item($tc2DataRow, 4).Value2
    $tc4Step = $ws.Cells.Item($tc4DataRow, 2).Value2
    $tc4Result = $ws.Cells.Item($tc4DataRow, 4).Value2

    $ws.Cells.Item($tc2DataRow, 2).Value = $tc4Step
    $ws.Cells.Item($tc2DataRow, 4).Value = $tc4Result
    $ws.Cells.Item($tc4DataRow, 2).Value = $tc2Step
    $ws.Cells.Item($tc4DataRow, 4).Value = $tc2Result
}
